$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the cells as text (matches the original numberStoredAsText content)
$ws.Range("A2:K2").NumberFormat = "@"

# Update row 2 with the data previously found in row 5
$ws.Range("A2").Value = " Nov 1 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "Super Kings won by 9 wickets (with 7 balls remaining)"
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Nicholas Pooran "
$ws.Range("G2").Value = "2"
$ws.Range("H2").Value = "6"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "33.33"

# Remove rows 3 through 6 (old data no longer needed)
$ws.Rows("3:6").Delete()
